$wb = $excel.ActiveWorkbook

# Insert the new "July Total" worksheet right before "July Report",
# mirroring the structure of the existing "June Total" summary sheet.
$julyReport = $wb.Worksheets.Item("July Report")
$julyTotal = $wb.Worksheets.Add($julyReport)
$julyTotal.Name = "July Total"

$julyTotal.Range("A1").Value = "Cash In Total"
$julyTotal.Range("B1").Value = "On-Us Check Total"
$julyTotal.Range("C1").Value = "Not On-Us Total"

$julyTotal.Range("A2").Value = 35630
$julyTotal.Range("B2").Value = 89075
$julyTotal.Range("C2").Value = 19565

$julyTotal.Activate()
